# Update cryptocurrency price/volume data to reflect latest values scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds values that can look numeric (e.g. "1.00", "9.00").
# Force it to Text format first so Excel does not silently convert these strings to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '71.208.36'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '3.683.75'
$ws.Range("E3").Value = '  +7.40%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '584.37'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = '179.23'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").Value = '3.671.53'
$ws.Range("E7").Value = '  +7.25%  '
$ws.Range("D8").Value = '0.612'
$ws.Range("E8").Value = '  +3.25%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("D11").Value = '0.610'
$ws.Range("E11").Value = '  +3.77%  '
$ws.Range("D12").Value = '49.39'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '0.0000285'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '4.279.11'
$ws.Range("E14").Value = '  +7.42%  '
$ws.Range("D15").Value = '681.53'
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("D16").Value = '9.00'
$ws.Range("E16").Value = '  +3.67%  '
$ws.Range("D17").Value = '3.682.96'
$ws.Range("E17").Value = '  +6.82%  '
$ws.Range("D18").Value = '71.439.62'
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").Value = '18.00'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '11.60'
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("D22").Value = '6.43'
$ws.Range("E22").Value = '  +19.74%  '
$ws.Range("D23").Value = '0.939'
$ws.Range("E23").Value = '  +3.85%  '
$ws.Range("D24").Value = '17.41'
$ws.Range("E24").Value = '  +1.68%  '
$ws.Range("D25").Value = '102.11'
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("D26").Value = '4.00'
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("D27").Value = '2.84'
$ws.Range("E27").Value = '  +5.24%  '
$ws.Range("D28").Value = '10.26'
$ws.Range("E28").Value = '  +5.65%  '
$ws.Range("D29").Value = '35.25'
$ws.Range("E29").Value = '  +4.39%  '
$ws.Range("D30").Value = '9.18'
$ws.Range("E30").Value = '  +4.11%  '
$ws.Range("D31").Value = '7.34'
$ws.Range("E31").Value = '  +4.85%  '
$ws.Range("D32").Value = '4.12'
$ws.Range("E32").Value = '  +8.28%  '
$ws.Range("D33").Value = '582.96'
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("D34").Value = '11.22'
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("D35").Value = '0.108'
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("D36").Value = '59.00'
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.147'
$ws.Range("E38").Value = '  +5.43%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.664.53'
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("D40").Value = '35.56'
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("E41").Value = '  +3.94%  '
$ws.Range("D42").Value = '3.46'
$ws.Range("E42").Value = '  +4.59%  '
$ws.Range("D43").Value = '0.0463'
$ws.Range("E43").Value = '  +8.88%  '
$ws.Range("E44").Value = '  +3.16%  '
$ws.Range("D45").Value = '0.348'
$ws.Range("E45").Value = '  +3.76%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").Value = '2.89'
$ws.Range("E46").Value = '  +8.25%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.39'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  +3.16%  '
$ws.Range("E49").Value = '  -3.24%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '137.02'
$ws.Range("E51").Value = '  +4.51%  '
